$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.542.56"
$ws.Range("E2").Value = "  -4.03%  "
$ws.Range("D3").Value = "3.321.01"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "572.75"
$ws.Range("E5").Value = "  -3.37%  "
$ws.Range("D6").Value = "182.59"
$ws.Range("E6").Value = "  -4.95%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.93%  "
$ws.Range("D9").Value = "0.129"
$ws.Range("E9").Value = "  -3.05%  "
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("E11").Value = "  -4.40%  "
$ws.Range("D12").Value = "3.899.54"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").Value = "27.17"
$ws.Range("E14").Value = "  -3.68%  "
$ws.Range("D15").Value = "66.611.56"
$ws.Range("E15").Value = "  -3.98%  "
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").Value = "3.301.16"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "438.51"
$ws.Range("E18").Value = "  +3.04%  "
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("E20").Value = "  -2.30%  "
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").Value = "73.91"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").Value = "3.461.92"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("D27").Value = "0.193"
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("D28").Value = "9.05"
$ws.Range("E28").Value = "  -5.65%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("D31").Value = "22.91"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").Value = "5.35"
$ws.Range("E32").Value = "  -4.66%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("E34").Value = "  -3.02%  "
$ws.Range("E35").Value = "  -3.93%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "160.19"
$ws.Range("E37").Value = "  -2.85%  "
$ws.Range("D38").Value = "27.52"
$ws.Range("E38").Value = "  +1.94%  "
$ws.Range("E39").Value = "  -3.93%  "
$ws.Range("D40").Value = "2.814.00"
$ws.Range("E41").Value = "  -2.12%  "
$ws.Range("E42").Value = "  -2.14%  "
$ws.Range("D43").Value = "6.24"
$ws.Range("E43").Value = "  -3.38%  "
$ws.Range("E44").Value = "  -0.93%  "
$ws.Range("D45").Value = "40.20"
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").Value = "24.24"
$ws.Range("E46").Value = "  -5.02%  "
$ws.Range("E47").Value = "  -6.85%  "
$ws.Range("D48").Value = "319.51"
$ws.Range("E48").Value = "  -6.93%  "
$ws.Range("E49").Value = "  -3.14%  "
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("D51").Value = "6.18"
$ws.Range("E51").Value = "  -1.60%  "
